$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.735.42'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '3.026.98'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.39'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '151.56'
$ws.Range('E6').Value = '  +6.82%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.024.85'
$ws.Range('E8').Value = '  +1.92%  '
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.41'
$ws.Range('E10').Value = '  +11.12%  '
$ws.Range('E11').Value = '  +5.73%  '
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.84'
$ws.Range('E14').Value = '  +2.67%  '
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('D16').Value = '3.529.25'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').Value = '62.684.16'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').Value = '3.034.46'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '450.41'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.19'
$ws.Range('E21').Value = '  +2.15%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.691'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('E23').Value = '  +2.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '82.37'
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('E25').Value = '  +5.36%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.88'
$ws.Range('E26').Value = '  +10.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.18'
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  +3.05%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.40'
$ws.Range('E30').Value = '  +8.04%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.15'
$ws.Range('E32').Value = '  +4.77%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.58'
$ws.Range('E33').Value = '  +2.65%  '
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('D35').Value = '0.0₃0863'
$ws.Range('E35').Value = '  +11.38%  '
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('E37').Value = '  +2.94%  '
$ws.Range('E38').Value = '  +12.60%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.09'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '50.36'
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('E42').Value = '  +4.43%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.294'
$ws.Range('E43').Value = '  +11.74%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '41.06'
$ws.Range('E44').Value = '  +10.80%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '392.66'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0355'
$ws.Range('D47').Value = '2.743.97'
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.37'
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('E51').Value = '  +0.46%  '
